# gsc-export/HTTPS.xlsx -- roll the "Chart" sheet's date window forward by one day.
#
# The export keeps a fixed-length trailing window of daily rows (A: Date,
# B: Non-HTTPS URLs, C: HTTPS URLs). Each new export drops the oldest date,
# shifts every remaining row up by one day, and appends a fresh row for the
# newest date. Column B (Non-HTTPS URLs) is always 0 and does not change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# New date labels for rows 2..91 (oldest "2025-10-07" dropped, "2026-01-05" appended)
$dates = @("2025-10-08", "2025-10-09", "2025-10-10", "2025-10-11", "2025-10-12", "2025-10-13", "2025-10-14", "2025-10-15", "2025-10-16", "2025-10-17", "2025-10-18", "2025-10-19", "2025-10-20", "2025-10-21", "2025-10-22", "2025-10-23", "2025-10-24", "2025-10-25", "2025-10-26", "2025-10-27", "2025-10-28", "2025-10-29", "2025-10-30", "2025-10-31", "2025-11-01", "2025-11-02", "2025-11-03", "2025-11-04", "2025-11-05", "2025-11-06", "2025-11-07", "2025-11-08", "2025-11-09", "2025-11-10", "2025-11-11", "2025-11-12", "2025-11-13", "2025-11-14", "2025-11-15", "2025-11-16", "2025-11-17", "2025-11-18", "2025-11-19", "2025-11-20", "2025-11-21", "2025-11-22", "2025-11-23", "2025-11-24", "2025-11-25", "2025-11-26", "2025-11-27", "2025-11-28", "2025-11-29", "2025-11-30", "2025-12-01", "2025-12-02", "2025-12-03", "2025-12-04", "2025-12-05", "2025-12-06", "2025-12-07", "2025-12-08", "2025-12-09", "2025-12-10", "2025-12-11", "2025-12-12", "2025-12-13", "2025-12-14", "2025-12-15", "2025-12-16", "2025-12-17", "2025-12-18", "2025-12-19", "2025-12-20", "2025-12-21", "2025-12-22", "2025-12-23", "2025-12-24", "2025-12-25", "2025-12-26", "2025-12-27", "2025-12-28", "2025-12-29", "2025-12-30", "2025-12-31", "2026-01-01", "2026-01-02", "2026-01-03", "2026-01-04", "2026-01-05")

# New "HTTPS URLs" counts for rows 2..91, shifted up by one row to match,
# with a new value for the newly appended "2026-01-05" row.
$values = @(1, 11, 16, 22, 33, 43, 51, 58, 58, 68, 66, 71, 71, 75, 82, 92, 84, 85, 86, 90, 83, 90, 93, 92, 102, 113, 115, 107, 105, 100, 94, 86, 83, 66, 54, 46, 43, 40, 37, 35, 30, 29, 26, 25, 25, 26, 26, 25, 25, 27, 28, 28, 27, 27, 27, 27, 27, 26, 25, 25, 25, 26, 27, 27, 29, 29, 30, 30, 31, 31, 31, 31, 31, 32, 32, 32, 32, 30, 31, 32, 30, 28, 28, 28, 28, 29, 29, 28, 27, 27)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    # Leading apostrophe forces the date-shaped string to stay text instead
    # of being auto-converted to a date serial by Excel's input parser.
    $ws.Cells.Item($row, 1).Value = "'" + $dates[$i]
    $ws.Cells.Item($row, 3).Value = $values[$i]
}

# Re-entering text through .Value applies a "Text" number format to column A;
# paste the original (General) format from column B back over it so the
# dates keep rendering exactly as they did before (no style drift).
$ws.Range("B2").Copy()
$ws.Range("A2:A91").PasteSpecial(-4122)
$excel.CutCopyMode = $false
